# Recalculate the "Median Value" ratio and corresponding "Tier" label for
# each school now that the score is computed relative to the median AFTER
# merging with the zip/census tract data (commit: "Update to calculate
# scores relative to median AFTER merging with zip/census tract").
#
# Column layout: A=School, B=District, C=Median Value, D=Tier

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updatedRows = @(
    @{ Row = 2; MedianValue = 1.028985507246377; Tier = "4th Tier" },
    @{ Row = 3; MedianValue = 0.5090579710144928; Tier = "Below Median" },
    @{ Row = 4; MedianValue = 1.957971014492754; Tier = "1st Tier" },
    @{ Row = 5; MedianValue = 1.318840579710145; Tier = "3rd Tier" },
    @{ Row = 6; MedianValue = 1.540760869565217; Tier = "1st Tier" },
    @{ Row = 7; MedianValue = 1; Tier = "4th Tier" },
    @{ Row = 8; MedianValue = 0.6898550724637681; Tier = "Below Median" },
    @{ Row = 9; MedianValue = 1.855676328502415; Tier = "1st Tier" },
    @{ Row = 10; MedianValue = 0.7608695652173912; Tier = "Below Median" },
    @{ Row = 11; MedianValue = 1.00054347826087; Tier = "4th Tier" },
    @{ Row = 12; MedianValue = 1.032608695652174; Tier = "4th Tier" },
    @{ Row = 13; MedianValue = 1.087409420289855; Tier = "4th Tier" },
    @{ Row = 14; MedianValue = 1.334692028985507; Tier = "2nd Tier" },
    @{ Row = 15; MedianValue = 1.458333333333333; Tier = "2nd Tier" },
    @{ Row = 16; MedianValue = 0.8293478260869566; Tier = "Below Median" },
    @{ Row = 17; MedianValue = 0.9146286231884058; Tier = "Below Median" },
    @{ Row = 18; MedianValue = 0.5727657004830917; Tier = "Below Median" },
    @{ Row = 19; MedianValue = 0.8510466988727858; Tier = "Below Median" },
    @{ Row = 20; MedianValue = 0.7059178743961352; Tier = "Below Median" },
    @{ Row = 21; MedianValue = 0.6518115942028985; Tier = "Below Median" },
    @{ Row = 22; MedianValue = 0.5217391304347826; Tier = "Below Median" },
    @{ Row = 23; MedianValue = 0.6105072463768116; Tier = "Below Median" },
    @{ Row = 24; MedianValue = 1.389855072463768; Tier = "2nd Tier" },
    @{ Row = 25; MedianValue = 1.675724637681159; Tier = "1st Tier" },
    @{ Row = 26; MedianValue = 1.361111111111111; Tier = "2nd Tier" },
    @{ Row = 27; MedianValue = 1.499547101449275; Tier = "1st Tier" },
    @{ Row = 28; MedianValue = 1.001811594202898; Tier = "4th Tier" },
    @{ Row = 29; MedianValue = 0.5757246376811593; Tier = "Below Median" },
    @{ Row = 30; MedianValue = 1.43677536231884; Tier = "2nd Tier" },
    @{ Row = 31; MedianValue = 1.059581320450886; Tier = "4th Tier" },
    @{ Row = 32; MedianValue = 1.678985507246377; Tier = "1st Tier" },
    @{ Row = 33; MedianValue = 1.016606280193237; Tier = "4th Tier" },
    @{ Row = 34; MedianValue = 0.9710144927536231; Tier = "Below Median" },
    @{ Row = 35; MedianValue = 0.4839975845410627; Tier = "Below Median" },
    @{ Row = 36; MedianValue = 1.356884057971014; Tier = "2nd Tier" },
    @{ Row = 37; MedianValue = 0.7355072463768115; Tier = "Below Median" },
    @{ Row = 38; MedianValue = 1.27536231884058; Tier = "3rd Tier" },
    @{ Row = 39; MedianValue = 1.306159420289855; Tier = "3rd Tier" },
    @{ Row = 40; MedianValue = 1.151449275362319; Tier = "3rd Tier" },
    @{ Row = 41; MedianValue = 1.109601449275362; Tier = "4th Tier" },
    @{ Row = 42; MedianValue = 0.8327294685990339; Tier = "Below Median" },
    @{ Row = 43; MedianValue = 1.534420289855072; Tier = "1st Tier" },
    @{ Row = 44; MedianValue = 0.8876811594202898; Tier = "Below Median" },
    @{ Row = 45; MedianValue = 1.113405797101449; Tier = "3rd Tier" },
    @{ Row = 46; MedianValue = 0.601086956521739; Tier = "Below Median" },
    @{ Row = 47; MedianValue = 0.9658816425120773; Tier = "Below Median" },
    @{ Row = 48; MedianValue = 1.185688405797101; Tier = "3rd Tier" },
    @{ Row = 49; MedianValue = 1.41268115942029; Tier = "2nd Tier" },
    @{ Row = 50; MedianValue = 1.071557971014493; Tier = "4th Tier" },
    @{ Row = 51; MedianValue = 0.8834541062801933; Tier = "Below Median" },
    @{ Row = 52; MedianValue = 0.6644927536231884; Tier = "Below Median" },
    @{ Row = 53; MedianValue = 1.204710144927536; Tier = "3rd Tier" },
    @{ Row = 54; MedianValue = 0.9035326086956521; Tier = "Below Median" },
    @{ Row = 55; MedianValue = 1.123188405797101; Tier = "3rd Tier" },
    @{ Row = 56; MedianValue = 0.9184782608695652; Tier = "Below Median" },
    @{ Row = 57; MedianValue = 0.527536231884058; Tier = "Below Median" },
    @{ Row = 58; MedianValue = 0.6124999999999999; Tier = "Below Median" },
    @{ Row = 59; MedianValue = 0.3043478260869565; Tier = "Below Median" },
    @{ Row = 60; MedianValue = 0.5126811594202898; Tier = "Below Median" },
    @{ Row = 61; MedianValue = 0.8211050724637681; Tier = "Below Median" },
    @{ Row = 62; MedianValue = 1.341032608695652; Tier = "2nd Tier" },
    @{ Row = 63; MedianValue = 0.6391304347826087; Tier = "Below Median" },
    @{ Row = 64; MedianValue = 0.5807971014492753; Tier = "Below Median" },
    @{ Row = 65; MedianValue = 0.6625905797101449; Tier = "Below Median" },
    @{ Row = 66; MedianValue = 0.3719806763285024; Tier = "Below Median" },
    @{ Row = 67; MedianValue = 0.9739130434782608; Tier = "Below Median" },
    @{ Row = 68; MedianValue = 1.911684782608696; Tier = "1st Tier" },
    @{ Row = 69; MedianValue = 0.6036231884057971; Tier = "Below Median" },
    @{ Row = 70; MedianValue = 0.8876811594202898; Tier = "Below Median" },
    @{ Row = 71; MedianValue = 1.607971014492753; Tier = "1st Tier" },
    @{ Row = 72; MedianValue = 1.557246376811594; Tier = "1st Tier" },
    @{ Row = 73; MedianValue = 0.5427536231884057; Tier = "Below Median" },
    @{ Row = 74; MedianValue = 1.327898550724637; Tier = "3rd Tier" },
    @{ Row = 75; MedianValue = 1.296014492753623; Tier = "3rd Tier" },
    @{ Row = 76; MedianValue = 1.348429951690821; Tier = "2nd Tier" },
    @{ Row = 77; MedianValue = 0.7684782608695652; Tier = "Below Median" },
    @{ Row = 78; MedianValue = 1.480525362318841; Tier = "2nd Tier" },
    @{ Row = 79; MedianValue = 0.6657608695652174; Tier = "Below Median" },
    @{ Row = 80; MedianValue = 0.6467391304347826; Tier = "Below Median" }
)

foreach ($row in $updatedRows) {
    $ws.Cells.Item($row.Row, 3).Value = $row.MedianValue
    $ws.Cells.Item($row.Row, 4).Value = $row.Tier
}
